$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.934.46'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.052.01'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '229.81'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.56'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +6.30%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.35%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.50%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.354.46'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.22%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.82'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.753'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.27'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.80%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.059.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.48%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.860.17'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.31%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.40'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.43'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.53%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.02%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.29'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '166.52'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +5.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.02'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.33'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.00%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.03%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.51%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +10.45%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.90%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.07'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +11.32%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.41%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.500.35'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.23'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.55'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0918'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.84%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +15.72%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.30%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.243.65'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.42%  '
